$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting (style) of column M (2019) onto the new column N (2020)
# so the new column inherits the same number format / font as its neighbour.
# Row 4 is a section header with no data in M, so it must stay untouched --
# skip it when mirroring formats from M to N.
$ws.Range("M3").Copy()
$ws.Range("N3").PasteSpecial(-4122)

$ws.Range("M5:M25").Copy()
$ws.Range("N5:N25").PasteSpecial(-4122)

# New 2020 data column
$ws.Cells.Item(3, 14).Value = 2020

$ws.Cells.Item(5, 14).Value = 2198.6999999999998
$ws.Cells.Item(6, 14).Value = 132.69999999999999
$ws.Cells.Item(7, 14).Value = 242.9
$ws.Cells.Item(8, 14).Value = 203.3
$ws.Cells.Item(9, 14).Value = 202.8
$ws.Cells.Item(10, 14).Value = 284.7
$ws.Cells.Item(11, 14).Value = 294.89999999999998
$ws.Cells.Item(12, 14).Value = 802.5
$ws.Cells.Item(13, 14).Value = 28.1
$ws.Cells.Item(14, 14).Value = 6.8

# Row 15 is an empty spacer row -- only the format is mirrored (done above),
# no value belongs there.

$ws.Cells.Item(16, 14).Value = 27.4
$ws.Cells.Item(17, 14).Value = 17.5
$ws.Cells.Item(18, 14).Value = 24.7
$ws.Cells.Item(19, 14).Value = 31.5
$ws.Cells.Item(20, 14).Value = 30.4
$ws.Cells.Item(21, 14).Value = 24.8
$ws.Cells.Item(22, 14).Value = 30.7
$ws.Cells.Item(23, 14).Value = 30.1
$ws.Cells.Item(24, 14).Value = 21.2
$ws.Cells.Item(25, 14).Value = 11.6

# Reflect the author's last selection (cell M25) before saving.
$ws.Range("M25").Select() | Out-Null
